$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the output values for the MobileNetV3(small) row
$ws.Range("E11").Value = 0.34549999999999997
$ws.Range("F11").Value = 0.88270000000000004

# Update the active selection to reflect where the user ended up
$ws.Range("F11").Select()
